$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove "Vai Trò" and "Trạng Thái" headers (L1, M1) - clear contents but keep style
$ws.Range("L1").ClearContents()
$ws.Range("M1").ClearContents()

# Remove "Nhân viên" and "Đang hoạt động" data values (L2, M2)
$ws.Range("L2").ClearContents()
$ws.Range("M2").ClearContents()

# Update selection to O3
$ws.Range("O3").Select()
